# Auto-generated: apply scheduled-runner price/profit refresh to the Leve profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3294.4707
$ws.Range("I76").Value = 3108.818
$ws.Range("J76").Value = 3634.8333
$ws.Range("K76").Value = 3108.818
$ws.Range("L76").Value = 3634.8333
$ws.Range("M76").Value = -2793.818
$ws.Range("N76").Value = -4264.8333

$ws.Range("H79").Value = 3294.4707
$ws.Range("I79").Value = 3108.818
$ws.Range("J79").Value = 3634.8333
$ws.Range("K79").Value = 3108.818
$ws.Range("L79").Value = 3634.8333
$ws.Range("M79").Value = -2016.818
$ws.Range("N79").Value = -5818.8333

$ws.Range("H112").Value = 3099999.2
$ws.Range("I112").Value = 1862.5
$ws.Range("J112").Value = 3439521.2
$ws.Range("K112").Value = 5587.5
$ws.Range("L112").Value = 10318563.6
$ws.Range("M112").Value = -4479.5
$ws.Range("N112").Value = -10320779.6

$ws.Range("H135").Value = 822.9
$ws.Range("I135").Value = 822.9
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7406.099999999999
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4871.099999999999
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 3336872.8
$ws.Range("I137").Value = 5559766.5
$ws.Range("J137").Value = 2532.25
$ws.Range("K137").Value = 16679299.5
$ws.Range("L137").Value = 7596.75
$ws.Range("M137").Value = -16676749.5
$ws.Range("N137").Value = -12696.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3129.1936
$ws.Range("I61").Value = 1845
$ws.Range("J61").Value = 4333.125
$ws.Range("K61").Value = 1845
$ws.Range("L61").Value = 4333.125
$ws.Range("M61").Value = -1633
$ws.Range("N61").Value = -4757.125

$ws.Range("H74").Value = 607.7917
$ws.Range("I74").Value = 592.2632
$ws.Range("J74").Value = 666.8
$ws.Range("K74").Value = 592.2632
$ws.Range("L74").Value = 666.8
$ws.Range("M74").Value = 281.7368
$ws.Range("N74").Value = -2414.8

$ws.Range("H77").Value = 607.7917
$ws.Range("I77").Value = 592.2632
$ws.Range("J77").Value = 666.8
$ws.Range("K77").Value = 2961.316
$ws.Range("L77").Value = 3334
$ws.Range("M77").Value = 1406.684
$ws.Range("N77").Value = -12070

$ws.Range("H132").Value = 23258208
$ws.Range("I132").Value = 27028862
$ws.Range("K132").Value = 81086586
$ws.Range("M132").Value = -81084056

$ws.Range("H136").Value = 3129.1936
$ws.Range("I136").Value = 1845
$ws.Range("J136").Value = 4333.125
$ws.Range("K136").Value = 5535
$ws.Range("L136").Value = 12999.375
$ws.Range("M136").Value = -2985
$ws.Range("N136").Value = -18099.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1980.2572
$ws.Range("I134").Value = 1687.2727
$ws.Range("J134").Value = 2476.077
$ws.Range("K134").Value = 5061.8181
$ws.Range("L134").Value = 7428.231000000001
$ws.Range("M134").Value = -2526.8181
$ws.Range("N134").Value = -12498.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1698513
$ws.Range("I31").Value = 2130244.2
$ws.Range("J31").Value = 7565.3335
$ws.Range("K31").Value = 2130244.2
$ws.Range("L31").Value = 7565.3335
$ws.Range("M31").Value = -2129949.2
$ws.Range("N31").Value = -8155.3335

$ws.Range("H34").Value = 1698513
$ws.Range("I34").Value = 2130244.2
$ws.Range("J34").Value = 7565.3335
$ws.Range("K34").Value = 2130244.2
$ws.Range("L34").Value = 7565.3335
$ws.Range("M34").Value = -2130042.2
$ws.Range("N34").Value = -7969.3335

$ws.Range("H58").Value = 8930999
$ws.Range("I58").Value = 1503.6666
$ws.Range("J58").Value = 25004092
$ws.Range("K58").Value = 1503.6666
$ws.Range("L58").Value = 25004092
$ws.Range("M58").Value = -1300.6666
$ws.Range("N58").Value = -25004498

$ws.Range("H132").Value = 2273.7346
$ws.Range("I132").Value = 1666.2122
$ws.Range("J132").Value = 3526.75
$ws.Range("K132").Value = 4998.6366
$ws.Range("L132").Value = 10580.25
$ws.Range("M132").Value = -2468.6366
$ws.Range("N132").Value = -15640.25

$ws.Range("H134").Value = 1164.1476
$ws.Range("I134").Value = 800.80554
$ws.Range("J134").Value = 1687.36
$ws.Range("K134").Value = 2402.41662
$ws.Range("L134").Value = 5062.08
$ws.Range("M134").Value = 132.58338
$ws.Range("N134").Value = -10132.08

$ws.Range("H136").Value = 8930999
$ws.Range("I136").Value = 1503.6666
$ws.Range("J136").Value = 25004092
$ws.Range("K136").Value = 4510.9998
$ws.Range("L136").Value = 75012276
$ws.Range("M136").Value = -1960.9998
$ws.Range("N136").Value = -75017376

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4610.8887
$ws.Range("I80").Value = 4583
$ws.Range("K80").Value = 4583
$ws.Range("M80").Value = -3585

$ws.Range("H83").Value = 4610.8887
$ws.Range("I83").Value = 4583
$ws.Range("K83").Value = 22915
$ws.Range("M83").Value = -17923

$ws.Range("H126").Value = 2554.0625
$ws.Range("I126").Value = 1488.0667
$ws.Range("J126").Value = 3494.647
$ws.Range("K126").Value = 4464.2001
$ws.Range("L126").Value = 10483.941
$ws.Range("M126").Value = -1994.2001
$ws.Range("N126").Value = -15423.941

$ws.Range("H132").Value = 2819.0227
$ws.Range("I132").Value = 2180.6553
$ws.Range("J132").Value = 4053.2
$ws.Range("K132").Value = 6541.965899999999
$ws.Range("L132").Value = 12159.6
$ws.Range("M132").Value = -4011.965899999999
$ws.Range("N132").Value = -17219.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 100001370
$ws.Range("J22").Value = 1837
$ws.Range("L22").Value = 1837
$ws.Range("N22").Value = -2427

$ws.Range("H27").Value = 100001370
$ws.Range("J27").Value = 1837
$ws.Range("L27").Value = 1837
$ws.Range("N27").Value = -2051

$ws.Range("H132").Value = 2359.6365
$ws.Range("I132").Value = 1337.8148
$ws.Range("K132").Value = 4013.4444
$ws.Range("M132").Value = -1483.4444

$ws.Range("H136").Value = 2441575.5
$ws.Range("I136").Value = 3847499.8
$ws.Range("J136").Value = 4640.2
$ws.Range("K136").Value = 11542499.4
$ws.Range("L136").Value = 13920.6
$ws.Range("M136").Value = -11539949.4
$ws.Range("N136").Value = -19020.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 174628.55
$ws.Range("I132").Value = 206343.53
$ws.Range("J132").Value = 45125.668
$ws.Range("K132").Value = 619030.59
$ws.Range("L132").Value = 135377.004
$ws.Range("M132").Value = -616500.59
$ws.Range("N132").Value = -140437.004

$ws.Range("H136").Value = 1167.0625
$ws.Range("I136").Value = 653.55554
$ws.Range("K136").Value = 1960.66662
$ws.Range("M136").Value = 589.33338
